$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 519.0470293638949
$ws.Range("D2").Value = 121.0308980644043
$ws.Range("G2").Value = 473
$ws.Range("H2").Value = 555
$ws.Range("C3").Value = 37.48204644382814
$ws.Range("D3").Value = 6.462295981724989
$ws.Range("F3").Value = 32.87
$ws.Range("G3").Value = 37.68
$ws.Range("H3").Value = 41.87
$ws.Range("C4").Value = 2.034764764627279
$ws.Range("D4").Value = 2.498994627991414
$ws.Range("H4").Value = 2.51
$ws.Range("C5").Value = 322.7068557850753
$ws.Range("D5").Value = 8.550966776585174
$ws.Range("F5").Value = 317.77
$ws.Range("G5").Value = 323.29
$ws.Range("H5").Value = 328.53
$ws.Range("C6").Value = 23.71782131545305
$ws.Range("D6").Value = 3.660608368792933
$ws.Range("F6").Value = 21.11
$ws.Range("G6").Value = 23.32
$ws.Range("H6").Value = 26.21
$ws.Range("C7").Value = -75.02723897846421
$ws.Range("D7").Value = 22.03234664286968
$ws.Range("F7").Value = -90
$ws.Range("G7").Value = -70
$ws.Range("C8").Value = 8.003956186816254
$ws.Range("D8").Value = 6.482616637954716
$ws.Range("F8").Value = 8.199999999999999
$ws.Range("H8").Value = 11.5
$ws.Range("C9").Value = 9.099572884200908
$ws.Range("D9").Value = 1.594939389211567
$ws.Range("C10").Value = 867.8223341045298
$ws.Range("D10").Value = 0.461020912626382
$ws.Range("C11").Value = 0.466567421990915
$ws.Range("D11").Value = 0.5292003112480498
$ws.Range("C12").Value = 22.75176248067642
$ws.Range("D12").Value = 12.29563296379348
$ws.Range("C13").Value = 0.6714732404741105
$ws.Range("D13").Value = 0.7500371934682567
$ws.Range("C14").Value = 1.831828711275575
$ws.Range("D14").Value = 1.669297308196148
$ws.Range("C15").Value = 92.28723897846365
$ws.Range("D15").Value = 22.03234664286967
$ws.Range("G15").Value = 87.25999999999999
$ws.Range("H15").Value = 107.26
$ws.Range("C16").Value = -84.46743065575667
$ws.Range("D16").Value = 19.90755133025647
$ws.Range("F16").Value = -100.265723755961
$ws.Range("G16").Value = -82.0778545523916
$ws.Range("H16").Value = -68.41392685158225
$ws.Range("C17").Value = -76.46347446894042
$ws.Range("D17").Value = 24.22605095149054
$ws.Range("F17").Value = -90.57382219273629
$ws.Range("G17").Value = -70.66683163887967
$ws.Range("H17").Value = -57.5149694202523
